$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from H1 into the new header cells I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data for columns I (I0) and J (IF), rows 2-54
$data = @(
    ,@(2, 5, 6)
    ,@(3, 7, 7)
    ,@(4, 4, 5)
    ,@(5, 7, 7)
    ,@(6, 4, 6)
    ,@(7, 7, 8)
    ,@(8, 5, 6)
    ,@(9, 7, 7)
    ,@(10, 7, 7)
    ,@(11, 11, 11)
    ,@(12, 5, 6)
    ,@(13, 5, 7)
    ,@(14, 5, 7)
    ,@(15, 8, 8)
    ,@(16, 3, 4)
    ,@(17, 8, 9)
    ,@(18, 5, 6)
    ,@(19, 6, 7)
    ,@(20, 9, 9)
    ,@(21, 7, 7)
    ,@(22, 5, 6)
    ,@(23, 3, 5)
    ,@(24, 7, 7)
    ,@(25, 8, 8)
    ,@(26, 7, 7)
    ,@(27, 5, 6)
    ,@(28, 8, 9)
    ,@(29, 5, 5)
    ,@(30, 7, 7)
    ,@(31, 4, 6)
    ,@(32, 6, 6)
    ,@(33, 6, 6)
    ,@(34, 10, 10)
    ,@(35, 9, 9)
    ,@(36, 2, 4)
    ,@(37, 8, 8)
    ,@(38, 2, 3)
    ,@(39, 5, 6)
    ,@(40, 9, 9)
    ,@(41, 6, 7)
    ,@(42, 1, 3)
    ,@(43, 9, 9)
    ,@(44, 8, 8)
    ,@(45, 8, 8)
    ,@(46, 8, 8)
    ,@(47, 7, 7)
    ,@(48, 3, 4)
    ,@(49, 9, 9)
    ,@(50, 8, 8)
    ,@(51, 5, 5)
    ,@(52, 8, 8)
    ,@(53, 6, 6)
    ,@(54, 6, 6)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}

